$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.243.24'
$ws.Range('E2').Value = '  +4.34%  '
$ws.Range('D3').Value = '3.609.27'
$ws.Range('E3').Value = '  +6.73%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '596.17'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('D6').Value = '183.89'
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('D7').Value = '3.599.87'
$ws.Range('E7').Value = '  +6.62%  '
$ws.Range('E8').Value = '  +2.11%  '
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E10').Value = '  +6.81%  '
$ws.Range('D11').Value = '0.608'
$ws.Range('E11').Value = '  +3.21%  '
$ws.Range('D12').Value = '50.18'
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('D13').Value = '0.0000291'
$ws.Range('E13').Value = '  +3.63%  '
$ws.Range('D14').Value = '696.53'
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').Value = '4.190.12'
$ws.Range('E15').Value = '  +6.89%  '
$ws.Range('D16').Value = '8.95'
$ws.Range('E16').Value = '  +3.97%  '
$ws.Range('D17').Value = '72.319.76'
$ws.Range('E17').Value = '  +4.46%  '
$ws.Range('D18').Value = '3.606.07'
$ws.Range('E18').Value = '  +6.41%  '
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').Value = '18.52'
$ws.Range('E20').Value = '  +5.03%  '
$ws.Range('E21').Value = '  +4.08%  '
$ws.Range('E22').Value = '  +3.47%  '
$ws.Range('E23').Value = '  +5.29%  '
$ws.Range('D24').Value = '17.69'
$ws.Range('E24').Value = '  +3.20%  '
$ws.Range('D25').Value = '105.11'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('D26').Value = '4.04'
$ws.Range('E26').Value = '  +2.76%  '
$ws.Range('E27').Value = '  +4.44%  '
$ws.Range('D28').Value = '10.15'
$ws.Range('E28').Value = '  +5.66%  '
$ws.Range('D29').Value = '35.18'
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('D30').Value = '9.07'
$ws.Range('E30').Value = '  +4.12%  '
$ws.Range('D31').Value = '7.49'
$ws.Range('E31').Value = '  +7.71%  '
$ws.Range('D32').Value = '4.16'
$ws.Range('E32').Value = '  +16.90%  '
$ws.Range('D33').Value = '595.88'
$ws.Range('E33').Value = '  +6.59%  '
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('D35').Value = '0.108'
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('D36').Value = '60.13'
$ws.Range('E36').Value = '  +2.50%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').Value = '3.665.93'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = '0.145'
$ws.Range('E39').Value = '  +4.89%  '
$ws.Range('E40').Value = '  +12.96%  '
$ws.Range('D41').Value = '36.20'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('E42').Value = '  +6.70%  '
$ws.Range('D43').Value = '2.83'
$ws.Range('E43').Value = '  +5.70%  '
$ws.Range('E44').Value = '  +4.51%  '
$ws.Range('E45').Value = '  +2.89%  '
$ws.Range('D46').Value = '3.38'
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('E47').Value = '  +4.25%  '
$ws.Range('D48').Value = '1.49'
$ws.Range('E48').Value = '  +5.51%  '
$ws.Range('E49').Value = '  +1.97%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').Value = '133.81'
$ws.Range('E51').Value = '  +0.56%  '
